$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 00:37"

$ws.Range("B4").Value = 5650855
$ws.Range("C4").Value = 38880
$ws.Range("D4").Value = 3001307
$ws.Range("E4").Value = 2474672
$ws.Range("G4").Value = 1160
$ws.Range("H4").Value = 174876

$ws.Range("D5").Value = 2554179
$ws.Range("E5").Value = 743287

$ws.Range("B11").Value = 489122
$ws.Range("C11").Value = 12462
$ws.Range("D11").Value = 312323
$ws.Range("E11").Value = 161180
$ws.Range("G11").Value = 247
$ws.Range("H11").Value = 15619

$ws.Range("B27").Value = 123071
$ws.Range("C27").Value = 199
$ws.Range("E27").Value = 4718
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 9045

$ws.Range("B28").Value = 115661
$ws.Range("C28").Value = 293
$ws.Range("D28").Value = 112355
$ws.Range("E28").Value = 3113

$ws.Range("B32").Value = 96753
$ws.Range("C32").Value = 163
$ws.Range("D32").Value = 61562
$ws.Range("E32").Value = 30007
$ws.Range("G32").Value = 11
$ws.Range("H32").Value = 5184

$ws.Range("B35").Value = 87123
$ws.Range("C35").Value = 386
$ws.Range("D35").Value = 55504
$ws.Range("E35").Value = 30130
$ws.Range("G35").Value = 8
$ws.Range("H35").Value = 1489

$ws.Range("B46").Value = 63847
$ws.Range("C46").Value = 903
$ws.Range("D46").Value = 52370
$ws.Range("E46").Value = 9058
$ws.Range("G46").Value = 30
$ws.Range("H46").Value = 2419

$ws.Range("B52").Value = 49895
$ws.Range("C52").Value = 410
$ws.Range("D52").Value = 37051
$ws.Range("E52").Value = 11863
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 981

$ws.Range("B81").Value = 14669
$ws.Range("C81").Value = 169
$ws.Range("D81").Value = 9699
$ws.Range("E81").Value = 4451
$ws.Range("G81").Value = 7
$ws.Range("H81").Value = 519

$ws.Range("B87").Value = 10111
$ws.Range("C87").Value = 51
$ws.Range("E87").Value = 992

$ws.Range("A92").Value = "Guinea"
$ws.Range("B92").Value = 8715
$ws.Range("C92").Value = 95
$ws.Range("D92").Value = 7532
$ws.Range("E92").Value = 1131
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 52

$ws.Range("A93").Value = "Guayana Francesa"
$ws.Range("B93").Value = 8657
$ws.Range("C93").Value = 35
$ws.Range("D93").Value = 8054
$ws.Range("E93").Value = 550
$ws.Range("H93").Value = 53

$ws.Range("B97").Value = 7921
$ws.Range("C97").Value = 24
$ws.Range("E97").Value = 2490

$ws.Range("A100").Value = "Luxemburgo"
$ws.Range("B100").Value = 7499
$ws.Range("C100").Value = 30
$ws.Range("D100").Value = 6753
$ws.Range("E100").Value = 622
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 124

$ws.Range("A101").Value = "Grecia"
$ws.Range("B101").Value = 7472
$ws.Range("C101").Value = 250
$ws.Range("D101").Value = 3804
$ws.Range("E101").Value = 3436
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 232

$ws.Range("B103").Value = 6789
$ws.Range("C103").Value = 27
$ws.Range("D103").Value = 6051
$ws.Range("E103").Value = 581

$ws.Range("B104").Value = 6079
$ws.Range("C104").Value = 170
$ws.Range("D104").Value = 3648
$ws.Range("E104").Value = 2407
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 24

$ws.Range("A105").Value = "Zimbabue"
$ws.Range("B105").Value = 5378
$ws.Range("C105").Value = 70
$ws.Range("D105").Value = 4105
$ws.Range("E105").Value = 1132
$ws.Range("G105").Value = 6
$ws.Range("H105").Value = 141

$ws.Range("A106").Value = "Republica de Yibuti"
$ws.Range("B106").Value = 5374
$ws.Range("C106").Value = 2
$ws.Range("D106").Value = 5216
$ws.Range("E106").Value = 99
$ws.Range("H106").Value = 59

$ws.Range("B115").Value = 3989
$ws.Range("C115").Value = 95
$ws.Range("D115").Value = 2587
$ws.Range("E115").Value = 1326
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 76

$ws.Range("A117").Value = "Cuba"
$ws.Range("B117").Value = 3408
$ws.Range("C117").Value = 44
$ws.Range("D117").Value = 2794
$ws.Range("E117").Value = 526
$ws.Range("H117").Value = 88

$ws.Range("A118").Value = "Tailandia"
$ws.Range("B118").Value = 3381
$ws.Range("C118").Value = 3
$ws.Range("D118").Value = 3198
$ws.Range("E118").Value = 125
$ws.Range("H118").Value = 58

$ws.Range("A151").Value = "Liberia"
$ws.Range("B151").Value = 1282
$ws.Range("C151").Value = 5
$ws.Range("D151").Value = 803
$ws.Range("E151").Value = 397
$ws.Range("H151").Value = 82

$ws.Range("A152").Value = "Burkina Faso"
$ws.Range("B152").Value = 1280
$ws.Range("D152").Value = 1018
$ws.Range("E152").Value = 207
$ws.Range("H152").Value = 55

$ws.Range("A154").Value = "Togo"
$ws.Range("B154").Value = 1173
$ws.Range("C154").Value = 19
$ws.Range("D154").Value = 868
$ws.Range("E154").Value = 278
$ws.Range("H154").Value = 27

$ws.Range("A155").Value = "Niger"
$ws.Range("B155").Value = 1167
$ws.Range("D155").Value = 1078
$ws.Range("E155").Value = 20
$ws.Range("H155").Value = 69

$ws.Range("A161").Value = "Reunion"
$ws.Range("B161").Value = 903
$ws.Range("C161").Value = 23
$ws.Range("D161").Value = 657
$ws.Range("E161").Value = 241
$ws.Range("H161").Value = 5

$ws.Range("A162").Value = "Santo Tome y Principe"
$ws.Range("B162").Value = 885
$ws.Range("D162").Value = 821
$ws.Range("E162").Value = 49
$ws.Range("H162").Value = 15

$ws.Range("A163").Value = "Guyana"
$ws.Range("B163").Value = 737
$ws.Range("C163").Value = 28
$ws.Range("D163").Value = 365
$ws.Range("E163").Value = 347
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 25

$ws.Range("A164").Value = "Crucero"
$ws.Range("B164").Value = 712
$ws.Range("D164").Value = 651
$ws.Range("E164").Value = 48
$ws.Range("H164").Value = 13

$ws.Range("B171").Value = 416
$ws.Range("C171").Value = 3
$ws.Range("E171").Value = 79

$ws.Range("B172").Value = 406
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 393
$ws.Range("E172").Value = 6

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
